$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9846819043159485
$ws.Range("B1").Value = 3.066795825958252
$ws.Range("C1").Value = 4.005835056304932
$ws.Range("D1").Value = 2.046697854995728
$ws.Range("E1").Value = 1.212541580200195
